$wb = $excel.ActiveWorkbook

# The "DeviceList" sheet has a column (E) for device
# "ONEPLUS_5T_Android_9.0.0_7e05c" that needs to be removed entirely.
$ws = $wb.Worksheets.Item("DeviceList")

# Shrink the conditional formatting range from B2:H2 to B2:G2 before the
# column shift happens (so it keeps referencing the correct relative cells).
$cfRange = $ws.Range("B2:H2")
$fc1 = $cfRange.FormatConditions.Item(1)
$fc2 = $cfRange.FormatConditions.Item(2)
$newCfRange = $ws.Range("B2:G2")
$fc1.ModifyAppliesToRange($newCfRange)
$fc2.ModifyAppliesToRange($newCfRange)

# Delete the whole column E, shifting F,G,H left to E,F,G.
$ws.Columns.Item(5).Delete()

# Move the active selection to match the new layout.
$ws.Range("B14").Select()
